$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values for rows 2-178 per the recorded run results.
# The remaining rows (179-252) already contain the correct value (7293) and are left untouched.

for ($r = 2; $r -le 3; $r++) {
    $ws.Cells.Item($r, 3).Value = 7343
}

$ws.Cells.Item(4, 3).Value = 7318

for ($r = 5; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = 7310
}

for ($r = 19; $r -le 178; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
